$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header E1 from "sourceRecordType" to "category"
$ws.Range("E1").Value = "category"

# Move the active selection from F4 to E6
$ws.Range("E6").Select()

# Nudge the sheet's default/standard column width
$ws.StandardWidth = 11.23046875

# F2 was a literal boolean TRUE; make it a live formula that evaluates to TRUE
$ws.Range("F2").Formula = "=TRUE()"
